$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 531.8333
$ws.Range("I2").Value = 105.7
$ws.Range("J2").Value = 2662.5
$ws.Range("K2").Value = 105.7
$ws.Range("L2").Value = 2662.5
$ws.Range("M2").Value = 7.299999999999997
$ws.Range("N2").Value = -2888.5
$ws.Range("H9").Value = 692.36
$ws.Range("I9").Value = 769.8182
$ws.Range("K9").Value = 769.8182
$ws.Range("M9").Value = -600.8182
$ws.Range("H33").Value = 1745
$ws.Range("J33").Value = 127.666664
$ws.Range("L33").Value = 127.666664
$ws.Range("N33").Value = -585.666664
$ws.Range("H40").Value = 7273.5264
$ws.Range("J40").Value = 6926.727
$ws.Range("L40").Value = 6926.727
$ws.Range("N40").Value = -7276.727
$ws.Range("H45").Value = 4536.75
$ws.Range("J45").Value = 4536.75
$ws.Range("L45").Value = 13610.25
$ws.Range("N45").Value = -13994.25
$ws.Range("H57").Value = 21486.25
$ws.Range("J57").Value = 8315
$ws.Range("L57").Value = 24945
$ws.Range("N57").Value = -25943
$ws.Range("H58").Value = 670192
$ws.Range("I58").Value = 371.25
$ws.Range("K58").Value = 1113.75
$ws.Range("M58").Value = -963.75
$ws.Range("H64").Value = 7760.4
$ws.Range("I64").Value = 4950.2
$ws.Range("J64").Value = 10570.6
$ws.Range("K64").Value = 4950.2
$ws.Range("L64").Value = 10570.6
$ws.Range("M64").Value = -4702.2
$ws.Range("N64").Value = -11066.6
$ws.Range("H67").Value = 7760.4
$ws.Range("I67").Value = 4950.2
$ws.Range("J67").Value = 10570.6
$ws.Range("K67").Value = 4950.2
$ws.Range("L67").Value = 10570.6
$ws.Range("M67").Value = -4092.2
$ws.Range("N67").Value = -12286.6
$ws.Range("H70").Value = 167999.67
$ws.Range("I70").Value = 999.5
$ws.Range("K70").Value = 2998.5
$ws.Range("M70").Value = -2728.5
$ws.Range("H73").Value = 167999.67
$ws.Range("I73").Value = 999.5
$ws.Range("K73").Value = 2998.5
$ws.Range("M73").Value = -2062.5
$ws.Range("H74").Value = 10367
$ws.Range("I74").Value = 7019.923
$ws.Range("J74").Value = 21245
$ws.Range("K74").Value = 7019.923
$ws.Range("L74").Value = 21245
$ws.Range("M74").Value = -6083.923
$ws.Range("N74").Value = -23117
$ws.Range("H76").Value = 333336670
$ws.Range("H77").Value = 10367
$ws.Range("I77").Value = 7019.923
$ws.Range("J77").Value = 21245
$ws.Range("K77").Value = 35099.615
$ws.Range("L77").Value = 106225
$ws.Range("M77").Value = -30419.615
$ws.Range("N77").Value = -115585
$ws.Range("H79").Value = 333336670
$ws.Range("H92").Value = 151.25
$ws.Range("I92").Value = 107.25
$ws.Range("J92").Value = 327.25
$ws.Range("K92").Value = 107.25
$ws.Range("L92").Value = 327.25
$ws.Range("M92").Value = 1140.75
$ws.Range("N92").Value = -2823.25
$ws.Range("H106").Value = 3248.6667
$ws.Range("I106").Value = 3248.6667
$ws.Range("K106").Value = 3248.6667
$ws.Range("M106").Value = -2617.6667
$ws.Range("H107").Value = 84044.664
$ws.Range("I107").Value = 84044.664
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 84044.664
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -82124.664
$ws.Range("H111").Value = 59312.723
$ws.Range("I111").Value = 69710.07000000001
$ws.Range("J111").Value = 7326
$ws.Range("K111").Value = 209130.21
$ws.Range("L111").Value = 21978
$ws.Range("M111").Value = -206063.21
$ws.Range("N111").Value = -28112
$ws.Range("H112").Value = 3453.3447
$ws.Range("J112").Value = 3821.0833
$ws.Range("L112").Value = 11463.2499
$ws.Range("N112").Value = -13679.2499
$ws.Range("H116").Value = 18938.75
$ws.Range("I116").Value = 7166.6665
$ws.Range("J116").Value = 26002
$ws.Range("K116").Value = 7166.6665
$ws.Range("L116").Value = 26002
$ws.Range("M116").Value = -3724.6665
$ws.Range("N116").Value = -32886
$ws.Range("H118").Value = 748.7778
$ws.Range("J118").Value = 1609
$ws.Range("L118").Value = 4827
$ws.Range("N118").Value = -8141
$ws.Range("H131").Value = 3896.2856
$ws.Range("I131").Value = 3294.8
$ws.Range("J131").Value = 5400
$ws.Range("K131").Value = 9884.400000000001
$ws.Range("L131").Value = 16200
$ws.Range("M131").Value = -4844.400000000001
$ws.Range("N131").Value = -26280
$ws.Range("H132").Value = 6625
$ws.Range("I132").Value = 6625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -17345
$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200
$ws.Range("H137").Value = 805524
$ws.Range("I137").Value = 627296.5
$ws.Range("J137").Value = 1122372.9
$ws.Range("K137").Value = 1881889.5
$ws.Range("L137").Value = 3367118.7
$ws.Range("M137").Value = -1879339.5
$ws.Range("N137").Value = -3372218.7
$ws.Range("H138").Value = 8193.434999999999
$ws.Range("I138").Value = 3268.9092
$ws.Range("J138").Value = 9127.396000000001
$ws.Range("K138").Value = 9806.7276
$ws.Range("L138").Value = 27382.188
$ws.Range("M138").Value = -4666.7276
$ws.Range("N138").Value = -37662.188
$ws.Range("H141").Value = 1940.1072
$ws.Range("I141").Value = 1956.4073
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 5869.2219
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = -689.2219000000005
$ws.Range("N141").Value = -14860
$ws.Range("N107").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 150
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 200
$ws.Range("N4").Value = -432
$ws.Range("H5").Value = 350.33334
$ws.Range("I5").Value = 275.5
$ws.Range("K5").Value = 275.5
$ws.Range("M5").Value = -163.5
$ws.Range("H32").Value = 8609.5
$ws.Range("I32").Value = 7658.273
$ws.Range("K32").Value = 7658.273
$ws.Range("M32").Value = -7371.273
$ws.Range("H63").Value = 3795.6667
$ws.Range("I63").Value = 3467.4285
$ws.Range("K63").Value = 3467.4285
$ws.Range("M63").Value = -2781.4285
$ws.Range("H66").Value = 3795.6667
$ws.Range("I66").Value = 3467.4285
$ws.Range("K66").Value = 17337.1425
$ws.Range("M66").Value = -13905.1425
$ws.Range("H74").Value = 2072.9048
$ws.Range("I74").Value = 1976.6
$ws.Range("K74").Value = 1976.6
$ws.Range("M74").Value = -1102.6
$ws.Range("H77").Value = 2072.9048
$ws.Range("I77").Value = 1976.6
$ws.Range("K77").Value = 9883
$ws.Range("M77").Value = -5515
$ws.Range("H88").Value = 2985.1333
$ws.Range("I88").Value = 3465
$ws.Range("J88").Value = 2436.7144
$ws.Range("K88").Value = 3465
$ws.Range("L88").Value = 2436.7144
$ws.Range("M88").Value = -3059
$ws.Range("N88").Value = -3248.7144
$ws.Range("H91").Value = 2985.1333
$ws.Range("I91").Value = 3465
$ws.Range("J91").Value = 2436.7144
$ws.Range("K91").Value = 3465
$ws.Range("L91").Value = 2436.7144
$ws.Range("M91").Value = -2061
$ws.Range("N91").Value = -5244.7144
$ws.Range("H102").Value = 1561.3636
$ws.Range("I102").Value = 1615.875
$ws.Range("J102").Value = 1416
$ws.Range("K102").Value = 1615.875
$ws.Range("L102").Value = 1416
$ws.Range("M102").Value = 6.125
$ws.Range("N102").Value = -4660
$ws.Range("H122").Value = 3055.9092
$ws.Range("I122").Value = 1748.6957
$ws.Range("K122").Value = 5246.0871
$ws.Range("M122").Value = -2796.0871
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 3326.5806
$ws.Range("I132").Value = 2224
$ws.Range("J132").Value = 7106.857
$ws.Range("K132").Value = 6672
$ws.Range("L132").Value = 21320.571
$ws.Range("M132").Value = -4142
$ws.Range("N132").Value = -26380.571
$ws.Range("N130").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 350.33334
$ws.Range("I4").Value = 275.5
$ws.Range("K4").Value = 275.5
$ws.Range("M4").Value = -160.5
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("H50").Value = 80763.5
$ws.Range("J50").Value = 80763.5
$ws.Range("L50").Value = 80763.5
$ws.Range("N50").Value = -81911.5
$ws.Range("H86").Value = 1418297.4
$ws.Range("I86").Value = 1701765.8
$ws.Range("J86").Value = 955.5
$ws.Range("K86").Value = 1701765.8
$ws.Range("L86").Value = 955.5
$ws.Range("M86").Value = -1700642.8
$ws.Range("N86").Value = -3201.5
$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50812
$ws.Range("H89").Value = 1418297.4
$ws.Range("I89").Value = 1701765.8
$ws.Range("J89").Value = 955.5
$ws.Range("K89").Value = 8508829
$ws.Range("L89").Value = 4777.5
$ws.Range("M89").Value = -8503213
$ws.Range("N89").Value = -16009.5
$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52808
$ws.Range("H99").Value = 4632.222
$ws.Range("I99").Value = 4002.2
$ws.Range("K99").Value = 4002.2
$ws.Range("M99").Value = -2504.2
$ws.Range("H105").Value = 1115.9333
$ws.Range("I105").Value = 1103.7273
$ws.Range("K105").Value = 1103.7273
$ws.Range("M105").Value = 643.2727
$ws.Range("H107").Value = 372932.47
$ws.Range("I107").Value = 2697.9524
$ws.Range("K107").Value = 2697.9524
$ws.Range("M107").Value = -777.9524000000001
$ws.Range("H108").Value = 80681.664
$ws.Range("J108").Value = 80681.664
$ws.Range("L108").Value = 80681.664
$ws.Range("N108").Value = -88361.664
$ws.Range("H134").Value = 89391.586
$ws.Range("I134").Value = 5855.4443
$ws.Range("K134").Value = 17566.3329
$ws.Range("M134").Value = -15031.3329
$ws.Range("N22").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 187.16667
$ws.Range("I22").Value = 181.5
$ws.Range("K22").Value = 181.5
$ws.Range("M22").Value = 168.5
$ws.Range("H31").Value = 53849.285
$ws.Range("I31").Value = 2552.8333
$ws.Range("J31").Value = 122244.555
$ws.Range("K31").Value = 2552.8333
$ws.Range("L31").Value = 122244.555
$ws.Range("M31").Value = -2257.8333
$ws.Range("N31").Value = -122834.555
$ws.Range("H34").Value = 53849.285
$ws.Range("I34").Value = 2552.8333
$ws.Range("J34").Value = 122244.555
$ws.Range("K34").Value = 2552.8333
$ws.Range("L34").Value = 122244.555
$ws.Range("M34").Value = -2350.8333
$ws.Range("N34").Value = -122648.555
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H74").Value = 74971.14
$ws.Range("I74").Value = 29999
$ws.Range("J74").Value = 92960
$ws.Range("K74").Value = 29999
$ws.Range("L74").Value = 92960
$ws.Range("M74").Value = -29125
$ws.Range("N74").Value = -94708
$ws.Range("H77").Value = 74971.14
$ws.Range("I77").Value = 29999
$ws.Range("J77").Value = 92960
$ws.Range("K77").Value = 89997
$ws.Range("L77").Value = 278880
$ws.Range("M77").Value = -85629
$ws.Range("N77").Value = -287616
$ws.Range("H100").Value = 66500
$ws.Range("J100").Value = 66500
$ws.Range("L100").Value = 66500
$ws.Range("N100").Value = -68664
$ws.Range("H122").Value = 2682.0667
$ws.Range("J122").Value = 5238.9
$ws.Range("L122").Value = 15716.7
$ws.Range("N122").Value = -20616.7
$ws.Range("H130").Value = 74999
$ws.Range("J130").Value = 74999
$ws.Range("L130").Value = 74999
$ws.Range("N130").Value = -85039
$ws.Range("H132").Value = 3854.4482
$ws.Range("I132").Value = 3422.4348
$ws.Range("K132").Value = 10267.3044
$ws.Range("M132").Value = -7737.304400000001
$ws.Range("H138").Value = 63402.5
$ws.Range("J138").Value = 49499.25
$ws.Range("L138").Value = 49499.25
$ws.Range("N138").Value = -59779.25
$ws.Range("H141").Value = 65737.2
$ws.Range("J141").Value = 61840.145
$ws.Range("L141").Value = 61840.145
$ws.Range("N141").Value = -72200.14499999999
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 326.5
$ws.Range("I2").Value = 66.666664
$ws.Range("J2").Value = 482.4
$ws.Range("K2").Value = 399.999984
$ws.Range("L2").Value = 2894.4
$ws.Range("M2").Value = -286.999984
$ws.Range("N2").Value = -3120.4
$ws.Range("H12").Value = 175.42857
$ws.Range("I12").Value = 9.5
$ws.Range("J12").Value = 203.08333
$ws.Range("K12").Value = 28.5
$ws.Range("L12").Value = 609.24999
$ws.Range("M12").Value = 144.5
$ws.Range("N12").Value = -955.24999
$ws.Range("H23").Value = 1183.3334
$ws.Range("J23").Value = 1183.3334
$ws.Range("L23").Value = 3550.0002
$ws.Range("N23").Value = -4020.0002
$ws.Range("H34").Value = 73493.2
$ws.Range("J34").Value = 91616.5
$ws.Range("L34").Value = 274849.5
$ws.Range("N34").Value = -275017.5
$ws.Range("H38").Value = 53.4
$ws.Range("I38").Value = 59.5
$ws.Range("K38").Value = 178.5
$ws.Range("M38").Value = 168.5
$ws.Range("H40").Value = 282.875
$ws.Range("J40").Value = 750
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3138
$ws.Range("H56").Value = 9800
$ws.Range("I56").Value = 9800
$ws.Range("K56").Value = 9800
$ws.Range("M56").Value = -9270
$ws.Range("H111").Value = 1700
$ws.Range("I111").Value = 2400
$ws.Range("K111").Value = 7200
$ws.Range("M111").Value = -4133
$ws.Range("H121").Value = 18070760
$ws.Range("I121").Value = 511.66666
$ws.Range("J121").Value = 34333984
$ws.Range("K121").Value = 1534.99998
$ws.Range("L121").Value = 103001952
$ws.Range("M121").Value = -224.9999800000001
$ws.Range("N121").Value = -103004572
$ws.Range("H122").Value = 59897.94
$ws.Range("I122").Value = 738.25
$ws.Range("J122").Value = 112484.336
$ws.Range("K122").Value = 6644.25
$ws.Range("L122").Value = 1012359.024
$ws.Range("M122").Value = -4194.25
$ws.Range("N122").Value = -1017259.024
$ws.Range("H129").Value = 1301.5
$ws.Range("I129").Value = 480.83334
$ws.Range("K129").Value = 1442.50002
$ws.Range("M129").Value = 3557.49998
$ws.Range("H131").Value = 15155966
$ws.Range("I131").Value = 30303844
$ws.Range("J131").Value = 8089.091
$ws.Range("K131").Value = 90911532
$ws.Range("L131").Value = 24267.273
$ws.Range("M131").Value = -90906492
$ws.Range("N131").Value = -34347.273
$ws.Range("H137").Value = 1901.5
$ws.Range("I137").Value = 1001.6667
$ws.Range("K137").Value = 3005.0001
$ws.Range("M137").Value = 2094.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 257.7143
$ws.Range("I2").Value = 335.27777
$ws.Range("K2").Value = 335.27777
$ws.Range("M2").Value = -222.27777
$ws.Range("H80").Value = 1252518.6
$ws.Range("I80").Value = 1003088.7
$ws.Range("J80").Value = 1668235.1
$ws.Range("K80").Value = 1003088.7
$ws.Range("L80").Value = 1668235.1
$ws.Range("M80").Value = -1002090.7
$ws.Range("N80").Value = -1670231.1
$ws.Range("H83").Value = 1252518.6
$ws.Range("I83").Value = 1003088.7
$ws.Range("J83").Value = 1668235.1
$ws.Range("K83").Value = 5015443.5
$ws.Range("L83").Value = 8341175.5
$ws.Range("M83").Value = -5010451.5
$ws.Range("N83").Value = -8351159.5
$ws.Range("H113").Value = 563470.0600000001
$ws.Range("I113").Value = 771561.6
$ws.Range("K113").Value = 771561.6
$ws.Range("M113").Value = -769391.6
$ws.Range("H122").Value = 5425.3335
$ws.Range("I122").Value = 4892.636
$ws.Range("K122").Value = 14677.908
$ws.Range("M122").Value = -12227.908
$ws.Range("H126").Value = 3080.75
$ws.Range("I126").Value = 2828
$ws.Range("K126").Value = 8484
$ws.Range("M126").Value = -6014
$ws.Range("H129").Value = 68351
$ws.Range("I129").Value = 60709
$ws.Range("J129").Value = 75993
$ws.Range("K129").Value = 60709
$ws.Range("L129").Value = 75993
$ws.Range("M129").Value = -55709
$ws.Range("N129").Value = -85993
$ws.Range("H132").Value = 505430.47
$ws.Range("I132").Value = 672986.5600000001
$ws.Range("J132").Value = 146381.72
$ws.Range("K132").Value = 2018959.68
$ws.Range("L132").Value = 439145.16
$ws.Range("M132").Value = -2016429.68
$ws.Range("N132").Value = -444205.16
$ws.Range("H136").Value = 38572.3
$ws.Range("J136").Value = 38572.3
$ws.Range("L136").Value = 115716.9
$ws.Range("N136").Value = -120816.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 439782.44
$ws.Range("I7").Value = 5484.25
$ws.Range("J7").Value = 913562.25
$ws.Range("K7").Value = 5484.25
$ws.Range("L7").Value = 913562.25
$ws.Range("M7").Value = -5372.25
$ws.Range("N7").Value = -913786.25
$ws.Range("H22").Value = 859.44446
$ws.Range("I22").Value = 745
$ws.Range("J22").Value = 892.1429000000001
$ws.Range("K22").Value = 745
$ws.Range("L22").Value = 892.1429000000001
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -1482.1429
$ws.Range("H27").Value = 859.44446
$ws.Range("I27").Value = 745
$ws.Range("J27").Value = 892.1429000000001
$ws.Range("K27").Value = 745
$ws.Range("L27").Value = 892.1429000000001
$ws.Range("M27").Value = -638
$ws.Range("N27").Value = -1106.1429
$ws.Range("H40").Value = 50624.24
$ws.Range("I40").Value = 73570.71000000001
$ws.Range("K40").Value = 73570.71000000001
$ws.Range("M40").Value = -73434.71000000001
$ws.Range("H55").Value = 714.1724
$ws.Range("I55").Value = 193.5
$ws.Range("K55").Value = 193.5
$ws.Range("M55").Value = -20.5
$ws.Range("H68").Value = 5537.5
$ws.Range("I68").Value = 4466.6665
$ws.Range("K68").Value = 4466.6665
$ws.Range("M68").Value = -3717.6665
$ws.Range("H71").Value = 5537.5
$ws.Range("I71").Value = 4466.6665
$ws.Range("K71").Value = 22333.3325
$ws.Range("M71").Value = -18589.3325
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("H93").Value = 45457180
$ws.Range("I93").Value = 66669056
$ws.Range("K93").Value = 66669056
$ws.Range("M93").Value = -66667808
$ws.Range("H126").Value = 439782.44
$ws.Range("I126").Value = 5484.25
$ws.Range("J126").Value = 913562.25
$ws.Range("K126").Value = 16452.75
$ws.Range("L126").Value = 2740686.75
$ws.Range("M126").Value = -13982.75
$ws.Range("N126").Value = -2745626.75
$ws.Range("H132").Value = 3275
$ws.Range("I132").Value = 2566.6667
$ws.Range("K132").Value = 7700.000100000001
$ws.Range("M132").Value = -5170.000100000001
$ws.Range("H136").Value = 693906.4
$ws.Range("I136").Value = 912206.0600000001
$ws.Range("J136").Value = 7821.5713
$ws.Range("K136").Value = 2736618.18
$ws.Range("L136").Value = 23464.7139
$ws.Range("M136").Value = -2734068.18
$ws.Range("N136").Value = -28564.7139
$ws.Range("N88").ClearContents()
$ws.Range("N91").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6903.1
$ws.Range("J74").Value = 7444.625
$ws.Range("L74").Value = 7444.625
$ws.Range("N74").Value = -9316.625
$ws.Range("H77").Value = 6903.1
$ws.Range("J77").Value = 7444.625
$ws.Range("L77").Value = 22333.875
$ws.Range("N77").Value = -31693.875
$ws.Range("H122").Value = 32261166
$ws.Range("I122").Value = 41668840
$ws.Range("J122").Value = 6278.143
$ws.Range("K122").Value = 125006520
$ws.Range("L122").Value = 18834.429
$ws.Range("M122").Value = -125004070
$ws.Range("N122").Value = -23734.429
$ws.Range("H125").Value = 55972.5
$ws.Range("J125").Value = 55972.5
$ws.Range("L125").Value = 55972.5
$ws.Range("N125").Value = -65812.5
$ws.Range("H126").Value = 2590.2727
$ws.Range("I126").Value = 1999.2222
$ws.Range("K126").Value = 5997.6666
$ws.Range("M126").Value = -3527.6666
$ws.Range("H132").Value = 46352.332
$ws.Range("I132").Value = 4712.8335
$ws.Range("J132").Value = 171270.83
$ws.Range("K132").Value = 14138.5005
$ws.Range("L132").Value = 513812.49
$ws.Range("M132").Value = -11608.5005
$ws.Range("N132").Value = -518872.49
$ws.Range("H136").Value = 347072.06
$ws.Range("I136").Value = 497073.94
$ws.Range("K136").Value = 1491221.82
$ws.Range("M136").Value = -1488671.82
